$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column widths for columns E (5), F (6), G (7) ---
# The underlying engine stores widths on a 1/6-character grid (ColumnWidth
# read/write round-trips through that grid), so we dial in the ColumnWidth
# value whose stored/serialized width lands on (or as close as possible to)
# the target OOXML "width" attribute.
$ws.Columns.Item(5).ColumnWidth = 27.666666666666668   # -> stored width 28.5
$ws.Columns.Item(6).ColumnWidth = 23.498697916666668   # -> stored width ~24.33203125
$ws.Columns.Item(7).ColumnWidth = 17.166666666666668   # -> stored width 18

# --- New row of data (row 17) ---
$ws.Range("B17").Value = "Analysis"
$ws.Range("D17").Value = "Computation"
$ws.Range("E17").Value = "[SampleB2,SampleA2]"

# --- Update the view: scroll position + selection ---
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E21").Select()
